$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.520.57'
$ws.Range("E2").Value = '  +1.06%  '

$ws.Range("D3").Value = '3.601.60'
$ws.Range("E3").Value = '  +1.90%  '

$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = $origStyle
$ws.Range("E4").Value = '  +0.06%  '

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '201.86'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  +3.80%  '

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '596.99'
$cell.Style = $origStyle
$ws.Range("E6").Value = '  -1.65%  '

$ws.Range("E7").Value = '  +0.68%  '

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = $origStyle
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  +6.84%  '

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.644'
$cell.Style = $origStyle
$ws.Range("E10").Value = '  -0.29%  '

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '53.90'
$cell.Style = $origStyle
$ws.Range("E11").Value = '  +1.14%  '

$ws.Range("E12").Value = '  -0.01%  '

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.64'
$cell.Style = $origStyle
$ws.Range("E13").Value = '  +1.90%  '

$ws.Range("D14").Value = '4.171.60'
$ws.Range("E14").Value = '  +1.99%  '

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '679.31'
$cell.Style = $origStyle
$ws.Range("E15").Value = '  +14.33%  '

$ws.Range("D16").Value = '70.637.47'
$ws.Range("E16").Value = '  +1.06%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.621.35'
$ws.Range("E17").Value = '  +3.83%  '

$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.15'
$cell.Style = $origStyle
$ws.Range("E18").Value = '  +0.87%  '

$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '12.79'
$cell.Style = $origStyle
$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("E20").Value = '  +0.56%  '

$ws.Range("E21").Value = '  +1.89%  '

$ws.Range("E22").Value = '  +5.85%  '

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '110.70'
$cell.Style = $origStyle
$ws.Range("E23").Value = '  +7.75%  '

$ws.Range("E24").Value = '  +3.52%  '

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.56'
$cell.Style = $origStyle
$ws.Range("E25").Value = '  -1.49%  '

$ws.Range("E26").Value = '  +0.06%  '

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.61'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  -1.24%  '

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.00'
$cell.Style = $origStyle
$ws.Range("E28").Value = '  -0.68%  '

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.12'
$cell.Style = $origStyle
$ws.Range("E29").Value = '  +6.33%  '

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '34.21'
$cell.Style = $origStyle
$ws.Range("E30").Value = '  +3.22%  '

$ws.Range("E31").Value = '  +6.66%  '

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.17'
$cell.Style = $origStyle
$ws.Range("E32").Value = '  +2.04%  '

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '12.33'
$cell.Style = $origStyle
$ws.Range("E33").Value = '  +0.21%  '

$ws.Range("E34").Value = '  +0.38%  '

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '63.58'
$cell.Style = $origStyle
$ws.Range("E35").Value = '  +0.26%  '

$ws.Range("B36").Value = 'PEPE'
$ws.Range("C36").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D36").Value = '0.0₃0853'
$ws.Range("E36").Value = '  +5.10%  '

$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '3.872.95'
$ws.Range("E37").Value = '  +1.49%  '

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range("E38").Value = '  -0.11%  '

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '512.93'
$cell.Style = $origStyle
$ws.Range("E39").Value = '  +0.49%  '

$ws.Range("E40").Value = '  -5.08%  '

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.59'
$cell.Style = $origStyle
$ws.Range("E41").Value = '  +0.77%  '

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '36.79'
$cell.Style = $origStyle
$ws.Range("E42").Value = '  +0.59%  '

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.385'
$cell.Style = $origStyle
$ws.Range("E43").Value = '  -1.47%  '

$ws.Range("E44").Value = '  +3.30%  '

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0468'
$cell.Style = $origStyle

$ws.Range("E46").Value = '  +9.70%  '

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.41'
$cell.Style = $origStyle
$ws.Range("E47").Value = '  +0.95%  '

$ws.Range("E48").Value = '  +1.76%  '

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.62'
$cell.Style = $origStyle
$ws.Range("E49").Value = '  +1.76%  '

$ws.Range("E50").Value = '  -0.23%  '

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.80'
$cell.Style = $origStyle
$ws.Range("E51").Value = '  +22.24%  '

Write-Output "Applied cryptos list update."